# Regenerate instance 00 and 01
# Updates three sheets: Productdata, Capacity, ProcessingTime

$wb = $excel.ActiveWorkbook

# --- Productdata sheet: C7, C8, C9 : 1 -> 2 ---
$wsProduct = $wb.Worksheets.Item("Productdata")
$wsProduct.Range("C7").Value = 2
$wsProduct.Range("C8").Value = 2
$wsProduct.Range("C9").Value = 2

# --- Capacity sheet: B3, B4, B6, B7, B9 ---
$wsCapacity = $wb.Worksheets.Item("Capacity")
$wsCapacity.Range("B3").Value = 40
$wsCapacity.Range("B4").Value = 100
$wsCapacity.Range("B6").Value = 80
$wsCapacity.Range("B7").Value = 5
$wsCapacity.Range("B9").Value = 5

# --- ProcessingTime sheet: C3, D4, F6, G7, I9 ---
$wsProcTime = $wb.Worksheets.Item("ProcessingTime")
$wsProcTime.Range("C3").Value = 4
$wsProcTime.Range("D4").Value = 5
$wsProcTime.Range("F6").Value = 4
$wsProcTime.Range("G7").Value = 1
$wsProcTime.Range("I9").Value = 1
